# The post "「クーフィー体」الكوفي" (originally row 33) was removed from the
# sheet. Deleting the entire row shifts every following row up by one,
# which is exactly what the target diff shows (old row 34 "「両手」..."
# becomes the new row 33, old row 35 becomes row 34, ..., old row 183
# becomes the new row 182), and updates the sheet dimension from
# A1:C183 to A1:C182.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(33).Delete()
